$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the count of "M" (male) users from 3 to 2
$ws.Range("B2").Value = 2

# Remove the last row (A4:B4, the "O" / 1 entry) entirely
$ws.Rows.Item(4).Delete()
